$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 6: swap the custom "Table_0" table style for the
#    built-in table style {900A5D2A-A7A9-436C-BDD4-E0D079A22558}.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tableShape = $s6.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{900A5D2A-A7A9-436C-BDD4-E0D079A22558}")

# ---------------------------------------------------------------------------
# 2) Theme colors: the deck's applied design ("Integral") is swapped for the
#    stock "Office Theme" palette - every slide picks this up because they
#    all share the one slide master / theme.
#    (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# ---------------------------------------------------------------------------
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$s1 = $p.Slides.Item(1)
$themeColors = $s1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgr = ($b * 65536) + ($g * 256) + $r
    $themeColors.Item($i).RGB = $bgr
}
